# Applies the cryptos.xlsx data refresh described in the commit diff.
# (GitHub Actions scheduled update of coin prices / 1h volume deltas)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.828.33'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.898.64'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = "'238.36"
$ws.Range('E5').Value = '  -3.05%  '
$ws.Range('D6').Value = "'0.9983"
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = "'0.4889"
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('D9').Value = "'0.06769"
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').Value = '1.907.41'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').Value = "'17.04"
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').Value = "'0.07284"
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'5.132"
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').Value = "'89.88"
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').Value = "'0.6696"
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').Value = '30.768.34'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = "'0.000007929"
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = "'0.9989"
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = '2.119.89'
$ws.Range('E20').Value = '  -0.74%  '
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.011"
$ws.Range('E22').Value = '  +3.09%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = "'209.46"
$ws.Range('E23').Value = '  +10.10%  '
$ws.Range('D24').Value = "'6.205"
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D25').Value = "'9.624"
$ws.Range('E25').Value = '  +2.87%  '
$ws.Range('D26').Value = "'158.50"
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('D27').Value = "'18.90"
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').Value = "'1.909"
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = "'1.421"
$ws.Range('E29').Value = '  +1.66%  '
$ws.Range('D30').Value = "'4.326"
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').Value = "'0.09153"
$ws.Range('E31').Value = '  +1.55%  '
$ws.Range('D32').Value = "'4.036"
$ws.Range('E32').Value = '  -0.39%  '
$ws.Range('D33').Value = "'0.05173"
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').Value = "'0.7497"
$ws.Range('E34').Value = '  +0.13%  '
$ws.Range('D35').Value = "'1.111"
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').Value = "'2.684"
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').Value = "'0.01836"
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('D38').Value = "'2.696"
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('D39').Value = "'0.9263"
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('D40').Value = "'2.097"
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('D41').Value = "'0.4504"
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('D42').Value = "'106.58"
$ws.Range('E42').Value = '  +0.40%  '
$ws.Range('D43').Value = "'5.802"
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('D44').Value = "'1.002"
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = "'7.807"
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').Value = "'0.1373"
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('D47').Value = "'66.88"
$ws.Range('E47').Value = '  +15.24%  '
$ws.Range('D48').Value = "'0.4087"
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = "'8.941"
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('D50').Value = "'0.05912"
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').Value = "'34.74"
$ws.Range('E51').Value = '  +3.48%  '
